$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A15").Value = "2025-04-23T09:15:35.407Z"
$ws.Range("B15").Value = "149.132.61.23"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = "Likely NOT Malignant"
